$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 - Interpretador do arquivo texto / Douglas
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 4).Value = "entregue dentro do prazo"
$ws.Cells.Item(2, 5).Value = "correção de bugs"

# ---------------------------------------------------------------------------
# Row 3 - Algoritmo de floyd warshall... / Tiago
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 4).Value = "entregue dentro do prazo"
# E3 keeps "Aguarda testes" (unchanged text)

# ---------------------------------------------------------------------------
# Row 4 - Algortimo de Dijkstra (menor caminho) / Ivens
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 4).Value = "entregue dentro do prazo"
$ws.Cells.Item(4, 5).Value = "finalizado"
$ws.Cells.Item(4, 5).HorizontalAlignment = -4108
$ws.Cells.Item(4, 5).VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 5 - Revisar documento de requisitos / Lilian
# ---------------------------------------------------------------------------
$ws.Cells.Item(5, 4).Value = ""
$ws.Cells.Item(5, 5).Value = "Aguardando início"

# ---------------------------------------------------------------------------
# Row 6 - Revisar diagrama de classes / Lilian
# ---------------------------------------------------------------------------
$ws.Cells.Item(6, 4).Value = ""
$ws.Cells.Item(6, 5).Value = "Aguardando início"

# ---------------------------------------------------------------------------
# Row 7 - Criar telas
# ---------------------------------------------------------------------------
$ws.Cells.Item(7, 2).Value = ""
$ws.Cells.Item(7, 4).Value = ""
$ws.Cells.Item(7, 5).Value = "Aguardando início"

# ---------------------------------------------------------------------------
# Row 8 - Artigo científico da aplicação
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 2).Value = ""
$ws.Cells.Item(8, 4).Value = ""
$ws.Cells.Item(8, 5).Value = "Aguardando início"

# ---------------------------------------------------------------------------
# Row 9 - Mostrar movimentação dos indivíduos (D9 keeps its text, untouched)
# ---------------------------------------------------------------------------
$ws.Cells.Item(9, 2).Value = ""

# ---------------------------------------------------------------------------
# Row 10 - now holds "Criar verificação..." + Douglas + delivery date
# ---------------------------------------------------------------------------
$ws.Cells.Item(10, 1).Value = "Criar verificação na leitura dos valores das arestas onde não pode ser negativa"
$ws.Cells.Item(10, 1).WrapText = $true
$ws.Cells.Item(10, 1).HorizontalAlignment = 1
$ws.Cells.Item(10, 1).VerticalAlignment = -4108
$ws.Cells.Item(10, 2).Value = "Douglas"
$ws.Cells.Item(10, 2).HorizontalAlignment = -4108
$ws.Cells.Item(10, 2).VerticalAlignment = -4108
$ws.Cells.Item(10, 3).Value = 43208
$ws.Cells.Item(10, 3).NumberFormat = "DD/MM/YY"
$ws.Cells.Item(10, 3).HorizontalAlignment = -4108
$ws.Cells.Item(10, 3).VerticalAlignment = -4108
$ws.Cells.Item(10, 4).Value = ""
$ws.Cells.Item(10, 5).Value = ""
$ws.Rows.Item(10).RowHeight = 35.95

# ---------------------------------------------------------------------------
# Row 11 - NEW task: Implementar escolha de direção da aresta / Ivens
# ---------------------------------------------------------------------------
$ws.Cells.Item(11, 1).Value = "Implementar escolha de direção da aresta"
$ws.Cells.Item(11, 2).Value = "Ivens"
$ws.Cells.Item(11, 2).HorizontalAlignment = -4108
$ws.Cells.Item(11, 2).VerticalAlignment = -4107
$ws.Cells.Item(11, 3).Value = 43208
$ws.Cells.Item(11, 3).NumberFormat = "DD/MM/YY"
$ws.Cells.Item(11, 3).HorizontalAlignment = -4108
$ws.Cells.Item(11, 3).VerticalAlignment = -4108
$ws.Rows.Item(11).RowHeight = 13.8

# ---------------------------------------------------------------------------
# Row 12 - NEW task: arquivos com grafos para teste / Tiago
# ---------------------------------------------------------------------------
$ws.Cells.Item(12, 1).Value = "arquivos com grafos para teste"
$ws.Cells.Item(12, 2).Value = "Tiago"
$ws.Cells.Item(12, 2).HorizontalAlignment = -4108
$ws.Cells.Item(12, 2).VerticalAlignment = -4107
$ws.Cells.Item(12, 3).Value = 43208
$ws.Cells.Item(12, 3).NumberFormat = "DD/MM/YY"
$ws.Cells.Item(12, 3).HorizontalAlignment = -4108
$ws.Cells.Item(12, 3).VerticalAlignment = -4107
$ws.Rows.Item(12).RowHeight = 13.8

# ---------------------------------------------------------------------------
# Column widths (approximate re-measure, mirroring the author's resave)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 60.7368421052632
$ws.Columns.Item(2).ColumnWidth = 16.2834008097166
$ws.Columns.Item(3).ColumnWidth = 24.2105263157895
$ws.Columns.Item(4).ColumnWidth = 31.0647773279352
$ws.Columns.Item(5).ColumnWidth = 22.6032388663968
$ws.Columns.Item(6).ColumnWidth = 17.4615384615385

# ---------------------------------------------------------------------------
# Selection, matching the diff's saved cursor position
# ---------------------------------------------------------------------------
$ws.Range("D6:D8").Select()
